$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numeric-looking Folder IDs that must stay text, not be
# auto-converted to numbers by Excel.
$ws.Range("A13:A49").NumberFormat = "@"

$data = @(
    @('21281662', 0.9997031464334705, 'acquisition_3 no existe', 0.9999998660010492, 'acquisition_3 no existe', 0.9997031464334705, 'acquisition_3 no existe'),
    @('22528025', 0.9885468955905812, 0.00006624206720765709, 0.999999866011401, 0.9980044446460157, 0.9775630411951303, 0.00006550604423868313),
    @('24259471', 0.9987777563443073, 0.0005762924382716049, 0.9999998658768962, 0.9997676039972112, 0.9987777563443073, 0.0005762924382716049),
    @('24780201', 0.9992728695130315, 0.0002176836848422497, 0.9999998659433507, 0.9993849938499385, 0.9992728695130315, 0.0002176836848422497),
    @('26444330', 0.9997602130486969, 0.000147355109739369, 0.999999866008698, 0.9990917347865577, 0.9997602130486969, 0.000147355109739369),
    @('26618355', 0.9988310721021948, 0.0002534507887517147, 0.9999998658840554, 0.9994717379820391, 0.9988310721021948, 0.0002534507887517147),
    @('26731283', 'acquisition_2 no existe', 'acquisition_3 no existe', 'acquisition_2 no existe', 'acquisition_3 no existe', 'acquisition_2 no existe', 'acquisition_3 no existe'),
    @('27040152', 0.9999232413837449, 0.00004166130829903978, 0.9999998660305439, 0.9967948717948718, 0.9999232413837449, 0.00004166130829903978),
    @('27130020', 0.9997363683127573, 'acquisition_3 no existe', 0.9999998660055021, 'acquisition_3 no existe', 0.9997363683127573, 'acquisition_3 no existe'),
    @('27434240', 'acquisition_2 no existe', 0.9999134623628257, 'acquisition_2 no existe', 0.9999998660292337, 'acquisition_2 no existe', 0.9999134623628257),
    @('27708269', 0.9995492273233882, 0.0001419967421124828, 0.9999998659804149, 0.9990574929311969, 0.9995492273233882, 0.0001419967421124828),
    @('27812007', 0.9998744802383402, 0.00004742155349794238, 0.9999998660240106, 0.9971830985915493, 0.9998744802383402, 0.00004742155349794238),
    @('27881548', 0.9998139306841564, 0.00007595486111111111, 0.999999866015897, 0.9982394366197183, 0.9998139306841564, 0.00007595486111111111),
    @('30114957', 0.9984563882458848, 0.0003331565072016461, 0.9999998658337267, 0.9995980707395499, 0.9984563882458848, 0.0003331565072016461),
    @('30345002', 0.9997722693758574, 0.00008171510631001372, 0.9999998660103138, 0.9983633387888707, 0.9997722693758574, 0.00008171510631001372),
    @('30385373', 0.9997663751714677, 'acquisition_3 no existe', 0.9999998660095237, 'acquisition_3 no existe', 0.9997663751714677, 'acquisition_3 no existe'),
    @('30401069', 0.9970834405006859, 0.001188084062071331, 0.9999998656489846, 0.9998872604284104, 0.9970834405006859, 0.001188084062071331),
    @('30860845', 0.998959271047668, 0.0006405928497942387, 0.9999998659012669, 0.9997909261969475, 0.998959271047668, 0.0006405928497942387),
    @('30974620', 0.999105956361454, 0.0003460165895061728, 0.9999998659209548, 0.9996130030959752, 0.999105956361454, 0.0003460165895061728),
    @('33068573', 0.9997498981910151, 0.00009712041323731139, 0.9999998660073155, 0.9986225895316805, 0.9997498981910151, 0.00009712041323731139),
    @('33212135', 0.9983938293038409, 0.0005513760288065844, 0.9999998658253199, 0.9997571046878795, 0.9983938293038409, 0.0005513760288065844),
    @('33275889', 0.9998645672582305, 'acquisition_3 no existe', 0.9999998660226823, 'acquisition_3 no existe', 0.9998645672582305, 'acquisition_3 no existe'),
    @('33327331', 0.9997912915809328, 'acquisition_3 no existe', 0.9999998660128631, 'acquisition_3 no existe', 0.9997912915809328, 'acquisition_3 no existe'),
    @('33604980', 0.9999018079132374, 0.0000369727366255144, 0.9999998660276722, 0.9963898916967509, 0.9999018079132374, 0.0000369727366255144),
    @('33637942', 0.9996681830847051, 0.00005746849279835391, 0.9999998659963626, 0.9976744186046511, 0.9996681830847051, 0.00005746849279835391),
    @('33773692', 0.9997216328017833, 0.0001332893947187929, 0.9999998660035271, 0.998995983935743, 0.9997216328017833, 0.0001332893947187929),
    @('33900783', 0.9997579357424554, 0.0000864036779835391, 0.9999998660083927, 0.9984520123839009, 0.9997579357424554, 0.0000864036779835391),
    @('33915934', 0.9995494952417695, 0.0001565982938957476, 0.9999998659804509, 0.9991452991452991, 0.9995494952417695, 0.0001565982938957476),
    @('34393044', 0.9996369705932785, 0.0001121238425925926, 0.9999998659921785, 0.9988066825775657, 0.9996369705932785, 0.0001121238425925926),
    @('34674493', 0.9996668434927983, 0.0001690564986282579, 0.9999998659961831, 0.9992082343626286, 0.9996668434927983, 0.0001690564986282579),
    @('35541276', 0.9984850555126886, 0.0004411276148834019, 0.9999998658375787, 0.9996964177292046, 0.9984850555126886, 0.0004411276148834019),
    @('36658002', 0.9986700531550069, 0.0002987289951989026, 0.9999998658624315, 0.9995517705064993, 0.9986700531550069, 0.0002987289951989026),
    @('36961015', 0.9997249817815501, 'acquisition_3 no existe', 0.999999866003976, 'acquisition_3 no existe', 0.9997249817815501, 'acquisition_3 no existe'),
    @('37252590', 0.9998877421982167, 'acquisition_3 no existe', 0.9999998660257876, 'acquisition_3 no existe', 0.9998877421982167, 'acquisition_3 no existe'),
    @('37467405', 'acquisition_2 no existe', 'acquisition_3 no existe', 'acquisition_2 no existe', 'acquisition_3 no existe', 'acquisition_2 no existe', 'acquisition_3 no existe'),
    @('37580279', 0.9999279299554183, 0.00001312800068587106, 0.9999998660311721, 0.98989898989899, 0.9999279299554183, 0.00001312800068587106),
    @('37764949', 0.9998625578703704, 'acquisition_3 no existe', 0.9999998660224131, 'acquisition_3 no existe', 0.9998625578703704, 'acquisition_3 no existe')
)

$r = 13
foreach ($row in $data) {
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

Write-Output "rows updated"